$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 459278, 683735850),
    @(3, 353, 517629),
    @(4, 516, 1131452),
    @(9, 1522, 4226770),
    @(11, 178519, 441484115),
    @(12, 327, 1350318),
    @(14, 92979, 219793018),
    @(15, 7, 24400),
    @(18, 5004, 8555008),
    @(21, 109, 312596),
    @(23, 11997, 24952950),
    @(25, 109165, 159647931),
    @(30, 480, 1234182),
    @(31, 48257, 115868070),
    @(34, 17529, 40647229),
    @(37, 2021, 3604221),
    @(39, 3076, 6088474),
    @(40, 134529, 200217214),
    @(42, 115, 229456),
    @(43, 5, 14868),
    @(46, 1699, 5172691),
    @(48, 67515, 167238661),
    @(51, 14381, 34546791),
    @(56, 3870, 7815439),
    @(57, 96797, 143906612),
    @(63, 626, 1714410),
    @(65, 41760, 98231139),
    @(68, 16792, 38123207),
    @(70, 1788, 3240506),
    @(74, 2546, 5151630),
    @(76, 28823, 44573604),
    @(80, 11942, 31821555),
    @(82, 8095, 19853112),
    @(84, 703, 1240215),
    @(85, 502, 996098),
    @(86, 201598, 305127766),
    @(87, 150, 473724),
    @(90, 777, 2112117),
    @(92, 94786, 224070573),
    @(93, 182, 735665),
    @(94, 28, 116085),
    @(95, 44853, 102697147),
    @(96, 23, 40459),
    @(97, 44, 280034),
    @(98, 7275, 27541302),
    @(100, 4749, 9634911),
    @(102, 44646, 66359480),
    @(106, 10852, 18880583),
    @(108, 10163, 16933245),
    @(110, 662, 1008635),
    @(112, 18648, 38146477),
    @(115, 4188, 9152634),
    @(117, 6078, 13721950),
    @(119, 271, 553097),
    @(120, 361, 779409),
    @(122, 201369, 293589909),
    @(128, 1630, 4558565),
    @(130, 78823, 187549044),
    @(131, 190, 684182),
    @(133, 42317, 96888140),
    @(136, 1835, 3740600),
    @(140, 3947, 8105217),
    @(142, 805159, 1286093512),
    @(143, 141, 252163),
    @(144, 347, 855743),
    @(147, 2814, 9840431),
    @(149, 326504, 788797478),
    @(150, 1097, 4809104),
    @(151, 63, 369850),
    @(152, 302922, 686201697),
    @(154, 66, 177144),
    @(155, 3858, 6514918),
    @(158, 12024, 25352020),
    @(161, 54980, 76486671),
    @(167, 16518, 26763496),
    @(169, 4573, 7250554),
    @(175, 24784, 37015452),
    @(176, 10, 12233),
    @(179, 10419, 21406803),
    @(181, 7327, 13894514),
    @(186, 33483, 81678710),
    @(188, 3263, 7653181),
    @(189, 435, 1011068),
    @(193, 124100, 186077918),
    @(197, 21, 30389),
    @(199, 1179, 3572341),
    @(201, 51516, 127033589),
    @(203, 20371, 48275232),
    @(205, 1598, 2768894),
    @(208, 2855, 5888660),
    @(210, 328601, 473817336),
    @(212, 241, 490223),
    @(217, 1452, 3796657),
    @(219, 127367, 299889658),
    @(220, 194, 687707),
    @(222, 50216, 114642338),
    @(225, 6614, 12099267),
    @(228, 8252, 15904822),
    @(231, 372820, 527844260),
    @(238, 985, 2647935),
    @(240, 141732, 332661823),
    @(243, 78535, 177957899),
    @(246, 6039, 10046699),
    @(249, 49, 197022),
    @(250, 10405, 20251612),
    @(253, 148506, 219513962),
    @(259, 1046, 3045586),
    @(261, 73752, 177153531),
    @(264, 19456, 44901083),
    @(266, 2397, 4386023),
    @(268, 4114, 8310781),
    @(269, 373825, 546796016),
    @(276, 1596, 5007866),
    @(278, 147839, 360532682),
    @(281, 104741, 244480512),
    @(284, 3309, 5745668),
    @(287, 8528, 17531275)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 3).Value = $u[1]
    $ws.Cells.Item($r, 4).Value = $u[2]
}
